$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: Full Name / EXAM 1 / EXAM 2 / EXAM 3 / HW TOTAL / QUIZ / PACT / TOTAL
$ws.Range("A1").Value = "Full Name"
$ws.Range("B1").Value = "EXAM 1"
$ws.Range("C1").Value = "EXAM 2"
$ws.Range("D1").Value = "EXAM 3"
$ws.Range("E1").Value = "HW TOTAL"
$ws.Range("F1").Value = "QUIZ"
$ws.Range("G1").Value = "PACT"
$ws.Range("H1").Value = "TOTAL"

# --- Per-student rows: B/C (exams 1-2) stay put, D is the new Exam 3 score,
#     E/F take over the old HW-total/Quiz values (shifted one column right),
#     G (PACT) is unchanged, H is the recomputed total out of 295.
$ws.Range("A2").Value = "Abrate, Tomas P."
$ws.Range("B2").Value = 66
$ws.Range("C2").Value = 95
$ws.Range("D2").Value = 84
$ws.Range("E2").Value = 74
$ws.Range("F2").Value = 15
$ws.Range("G2").Value = 5
$ws.Range("H2").Formula = "=ROUND(SUM(B2:G2)/295 * 100, 0)"

$ws.Range("A3").Value = "Angulo, Douglas J."
$ws.Range("B3").Value = 98
$ws.Range("C3").Value = 78
$ws.Range("D3").Value = 99
$ws.Range("E3").Value = 79
$ws.Range("F3").Value = 15
$ws.Range("G3").Value = 5
$ws.Range("H3").Formula = "=ROUND(SUM(B3:G3)/295 * 100, 0)"

$ws.Range("A4").Value = "Bannon, Mackin J."
$ws.Range("B4").Value = 99
$ws.Range("C4").Value = 98
$ws.Range("D4").Value = 83
$ws.Range("E4").Value = 78
$ws.Range("F4").Value = 15
$ws.Range("G4").Value = 5
$ws.Range("H4").Formula = "=ROUND(SUM(B4:G4)/295 * 100, 0)"

$ws.Range("A5").Value = "Barr, Michael A."
$ws.Range("B5").Value = 86
$ws.Range("C5").Value = 93
$ws.Range("D5").Value = 94
$ws.Range("E5").Value = 78
$ws.Range("F5").Value = 15
$ws.Range("G5").Value = 5
$ws.Range("H5").Formula = "=ROUND(SUM(B5:G5)/295 * 100, 0)"

$ws.Range("A6").Value = "Buckley, Emma E."
$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 90
$ws.Range("D6").Value = 93
$ws.Range("E6").Value = 75
$ws.Range("F6").Value = 15
$ws.Range("G6").Value = 5
$ws.Range("H6").Formula = "=ROUND(SUM(B6:G6)/295 * 100, 0)"

$ws.Range("A7").Value = "Darche, Michael K."
$ws.Range("B7").Value = 89
$ws.Range("C7").Value = 72
$ws.Range("D7").Value = 85
$ws.Range("E7").Value = 74
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = 5
$ws.Range("H7").Formula = "=ROUND(SUM(B7:G7)/295 * 100, 0)"

$ws.Range("A8").Value = "Dotzel, Sarah C."
$ws.Range("B8").Value = 100
$ws.Range("C8").Value = 77
$ws.Range("D8").Value = 89
$ws.Range("E8").Value = 78
$ws.Range("F8").Value = 15
$ws.Range("G8").Value = 5
$ws.Range("H8").Formula = "=ROUND(SUM(B8:G8)/295 * 100, 0)"

$ws.Range("A9").Value = "Faubert, Jonathan"
$ws.Range("B9").Value = 100
$ws.Range("C9").Value = 85
$ws.Range("D9").Value = 100
$ws.Range("E9").Value = 80
$ws.Range("F9").Value = 15
$ws.Range("G9").Value = 5
$ws.Range("H9").Formula = "=ROUND(SUM(B9:G9)/295 * 100, 0)"

$ws.Range("A10").Value = "Fazio, Louis L."
$ws.Range("B10").Value = 78
$ws.Range("C10").Value = 72
$ws.Range("D10").Value = 83
$ws.Range("E10").Value = 74
$ws.Range("F10").Value = 15
$ws.Range("G10").Value = 5
$ws.Range("H10").Formula = "=ROUND(SUM(B10:G10)/295 * 100, 0)"

$ws.Range("A11").Value = "Ganshirt, Lauren E."
$ws.Range("B11").Value = 91
$ws.Range("C11").Value = 80
$ws.Range("D11").Value = 95
$ws.Range("E11").Value = 79
$ws.Range("F11").Value = 15
$ws.Range("G11").Value = 5
$ws.Range("H11").Formula = "=ROUND(SUM(B11:G11)/295 * 100, 0)"

$ws.Range("A12").Value = "Grady, Martin J."
$ws.Range("B12").Value = 87
$ws.Range("C12").Value = 64
$ws.Range("D12").Value = 84
$ws.Range("E12").Value = 76
$ws.Range("F12").Value = 15
$ws.Range("G12").Value = 5
$ws.Range("H12").Formula = "=ROUND(SUM(B12:G12)/295 * 100, 0)"

$ws.Range("A13").Value = "Gu, Sirui"
$ws.Range("B13").Value = 93
$ws.Range("C13").Value = 69
$ws.Range("D13").Value = 91
$ws.Range("E13").Value = 69
$ws.Range("F13").Value = 15
$ws.Range("G13").Value = 5
$ws.Range("H13").Formula = "=ROUND(SUM(B13:G13)/295 * 100, 0)"

$ws.Range("A14").Value = "Gutierrez, Cristina"
$ws.Range("B14").Value = 99
$ws.Range("C14").Value = 83
$ws.Range("D14").Value = 99
$ws.Range("E14").Value = 79
$ws.Range("F14").Value = 15
$ws.Range("G14").Value = 5
$ws.Range("H14").Formula = "=ROUND(SUM(B14:G14)/295 * 100, 0)"

$ws.Range("A15").Value = "Healy, Conor P."
$ws.Range("B15").Value = 79
$ws.Range("C15").Value = 50
$ws.Range("D15").Value = 93
$ws.Range("E15").Value = 77
$ws.Range("F15").Value = 15
$ws.Range("G15").Value = 5
$ws.Range("H15").Formula = "=ROUND(SUM(B15:G15)/295 * 100, 0)"

$ws.Range("A16").Value = "Helme, Joseph W."
$ws.Range("B16").Value = 78
$ws.Range("C16").Value = 71
$ws.Range("D16").Value = 52
$ws.Range("E16").Value = 52
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 5
$ws.Range("H16").Formula = "=ROUND(SUM(B16:G16)/295 * 100, 0)"

$ws.Range("A17").Value = "Holdren, Kalie A."
$ws.Range("B17").Value = 71
$ws.Range("C17").Value = 68
$ws.Range("D17").Value = 80
$ws.Range("E17").Value = 75
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 5
$ws.Range("H17").Formula = "=ROUND(SUM(B17:G17)/295 * 100, 0)"

$ws.Range("A18").Value = "Howard, Kelsey R."
$ws.Range("B18").Value = 98
$ws.Range("C18").Value = 77
$ws.Range("D18").Value = 98
$ws.Range("E18").Value = 79
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 5
$ws.Range("H18").Formula = "=ROUND(SUM(B18:G18)/295 * 100, 0)"

$ws.Range("A19").Value = "Jones, Maura E."
$ws.Range("B19").Value = 80
$ws.Range("C19").Value = 73
$ws.Range("D19").Value = 74
$ws.Range("E19").Value = 73
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 5
$ws.Range("H19").Formula = "=ROUND(SUM(B19:G19)/295 * 100, 0)"

$ws.Range("A20").Value = "Kaes, Emily A."
$ws.Range("B20").Value = 76
$ws.Range("C20").Value = 71
$ws.Range("D20").Value = 83
$ws.Range("E20").Value = 66
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 5
$ws.Range("H20").Formula = "=ROUND(SUM(B20:G20)/295 * 100, 0)"

$ws.Range("A21").Value = "Kopetsky, Emily G."
$ws.Range("B21").Value = 96
$ws.Range("C21").Value = 68
$ws.Range("D21").Value = 97
$ws.Range("E21").Value = 78
$ws.Range("F21").Value = 15
$ws.Range("G21").Value = 5
$ws.Range("H21").Formula = "=ROUND(SUM(B21:G21)/295 * 100, 0)"

$ws.Range("A22").Value = "Larme, Marye C."
$ws.Range("B22").Value = 96
$ws.Range("C22").Value = 72
$ws.Range("D22").Value = 70
$ws.Range("E22").Value = 75
$ws.Range("F22").Value = 15
$ws.Range("G22").Value = 5
$ws.Range("H22").Formula = "=ROUND(SUM(B22:G22)/295 * 100, 0)"

$ws.Range("A23").Value = "Lattal, Sarah L."
$ws.Range("B23").Value = 83
$ws.Range("C23").Value = 71
$ws.Range("D23").Value = 80
$ws.Range("E23").Value = 75
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 5
$ws.Range("H23").Formula = "=ROUND(SUM(B23:G23)/295 * 100, 0)"

$ws.Range("A24").Value = "Orr, Justin M."
$ws.Range("B24").Value = 90
$ws.Range("C24").Value = 74
$ws.Range("D24").Value = 94
$ws.Range("E24").Value = 74
$ws.Range("F24").Value = 15
$ws.Range("G24").Value = 5
$ws.Range("H24").Formula = "=ROUND(SUM(B24:G24)/295 * 100, 0)"

$ws.Range("A25").Value = "Rossi, Michael N."
$ws.Range("B25").Value = 92
$ws.Range("C25").Value = 43
$ws.Range("D25").Value = 72
$ws.Range("E25").Value = 64
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 5
$ws.Range("H25").Formula = "=ROUND(SUM(B25:G25)/295 * 100, 0)"

$ws.Range("A26").Value = "Schoenwetter, Anne L."
$ws.Range("B26").Value = 95
$ws.Range("C26").Value = 78
$ws.Range("D26").Value = 93
$ws.Range("E26").Value = 76
$ws.Range("F26").Value = 15
$ws.Range("G26").Value = 5
$ws.Range("H26").Formula = "=ROUND(SUM(B26:G26)/295 * 100, 0)"

$ws.Range("A27").Value = "Smith, Allison M."
$ws.Range("B27").Value = 68
$ws.Range("C27").Value = 83
$ws.Range("D27").Value = 69
$ws.Range("E27").Value = 77
$ws.Range("F27").Value = 15
$ws.Range("G27").Value = 5
$ws.Range("H27").Formula = "=ROUND(SUM(B27:G27)/295 * 100, 0)"

$ws.Range("A28").Value = "Smith, Ellen D."
$ws.Range("B28").Value = 95
$ws.Range("C28").Value = 90
$ws.Range("D28").Value = 98
$ws.Range("E28").Value = 75
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 5
$ws.Range("H28").Formula = "=ROUND(SUM(B28:G28)/295 * 100, 0)"

$ws.Range("A29").Value = "Thomas, Grace M."
$ws.Range("B29").Value = 74
$ws.Range("C29").Value = 80
$ws.Range("D29").Value = 83
$ws.Range("E29").Value = 77
$ws.Range("F29").Value = 15
$ws.Range("G29").Value = 5
$ws.Range("H29").Formula = "=ROUND(SUM(B29:G29)/295 * 100, 0)"

$ws.Range("A30").Value = "Toohill, Connor P."
$ws.Range("B30").Value = 87
$ws.Range("C30").Value = 83
$ws.Range("D30").Value = 83
$ws.Range("E30").Value = 65
$ws.Range("F30").Value = 15
$ws.Range("G30").Value = 5
$ws.Range("H30").Formula = "=ROUND(SUM(B30:G30)/295 * 100, 0)"

$ws.Range("A31").Value = "Uber, Matthew D."
$ws.Range("B31").Value = 81
$ws.Range("C31").Value = 74
$ws.Range("D31").Value = 71
$ws.Range("E31").Value = 72
$ws.Range("F31").Value = 15
$ws.Range("G31").Value = 5
$ws.Range("H31").Formula = "=ROUND(SUM(B31:G31)/295 * 100, 0)"

$ws.Range("A32").Value = "Voutsos, Thomas L."
$ws.Range("B32").Value = 100
$ws.Range("C32").Value = 86
$ws.Range("D32").Value = 99
$ws.Range("E32").Value = 76
$ws.Range("F32").Value = 15
$ws.Range("G32").Value = 5
$ws.Range("H32").Formula = "=ROUND(SUM(B32:G32)/295 * 100, 0)"

$ws.Range("A33").Value = "Wasik, Peter M."
$ws.Range("B33").Value = 83
$ws.Range("C33").Value = 100
$ws.Range("D33").Value = 88
$ws.Range("E33").Value = 78
$ws.Range("F33").Value = 15
$ws.Range("G33").Value = 5
$ws.Range("H33").Formula = "=ROUND(SUM(B33:G33)/295 * 100, 0)"

$ws.Range("A34").Value = "Whichard, Johnny V."
$ws.Range("B34").Value = 70
$ws.Range("C34").Value = 71
$ws.Range("D34").Value = 76
$ws.Range("E34").Value = 78
$ws.Range("F34").Value = 15
$ws.Range("G34").Value = 5
$ws.Range("H34").Formula = "=ROUND(SUM(B34:G34)/295 * 100, 0)"

$ws.Range("A35").Value = "Yoo, Dong Suk"
$ws.Range("B35").Value = 95
$ws.Range("C35").Value = 77
$ws.Range("D35").Value = 97
$ws.Range("E35").Value = 76
$ws.Range("F35").Value = 15
$ws.Range("G35").Value = 5
$ws.Range("H35").Formula = "=ROUND(SUM(B35:G35)/295 * 100, 0)"

# --- Median / Mean summary rows, now spanning A:H
$ws.Range("A37").Value = "Median"
$ws.Range("B37").Formula = "=MEDIAN(B2:B35)"
$ws.Range("C37").Formula = "=MEDIAN(C2:C35)"
$ws.Range("D37").Formula = "=MEDIAN(D2:D35)"
$ws.Range("E37").Formula = "=MEDIAN(E2:E35)"
$ws.Range("F37").Formula = "=MEDIAN(F2:F35)"
$ws.Range("G37").Formula = "=MEDIAN(G2:G35)"
$ws.Range("H37").Formula = "=MEDIAN(H2:H35)"

$ws.Range("A38").Value = "Mean"
$ws.Range("B38").Formula = "=ROUNDUP(AVERAGE(B2:B35),1)"
$ws.Range("C38").Formula = "=ROUNDUP(AVERAGE(C2:C35),1)"
$ws.Range("D38").Formula = "=ROUNDUP(AVERAGE(D2:D35),1)"
$ws.Range("E38").Formula = "=ROUNDUP(AVERAGE(E2:E35),1)"
$ws.Range("F38").Formula = "=ROUNDUP(AVERAGE(F2:F35),1)"
$ws.Range("G38").Formula = "=ROUNDUP(AVERAGE(G2:G35),1)"
$ws.Range("H38").Formula = "=ROUNDUP(AVERAGE(H2:H35),1)"

# --- Cosmetic sheet-view tweaks from the diff: narrower PACT column, moved selection
$ws.Columns("G:G").ColumnWidth = 7.71
$ws.Range("E9").Select()

$wb.Save()
